$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1974.8334
$ws.Range("I100").Value = 1974.8334
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1974.8334
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1433.8334
$ws.Range("N100").ClearContents()

$ws.Range("H116").Value = 15006.2
$ws.Range("I116").Value = 51500
$ws.Range("J116").Value = 5882.75
$ws.Range("K116").Value = 51500
$ws.Range("L116").Value = 5882.75
$ws.Range("M116").Value = -48058
$ws.Range("N116").Value = -12766.75

$ws.Range("H129").Value = 864.0755
$ws.Range("I129").Value = 765.8333
$ws.Range("J129").Value = 876.617
$ws.Range("K129").Value = 2297.4999
$ws.Range("L129").Value = 2629.851
$ws.Range("M129").Value = 2702.5001
$ws.Range("N129").Value = -12629.851

$ws.Range("H132").Value = 1119.4222
$ws.Range("I132").Value = 950.46344
$ws.Range("J132").Value = 2851.25
$ws.Range("K132").Value = 2851.39032
$ws.Range("L132").Value = 8553.75
$ws.Range("M132").Value = -321.39032
$ws.Range("N132").Value = -13613.75

$ws.Range("H138").Value = 4036.476
$ws.Range("I138").Value = 3864.4
$ws.Range("J138").Value = 4466.6665
$ws.Range("K138").Value = 11593.2
$ws.Range("L138").Value = 13399.9995
$ws.Range("M138").Value = -6453.200000000001
$ws.Range("N138").Value = -23679.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1341.5143
$ws.Range("I74").Value = 1111.7931
$ws.Range("J74").Value = 2451.8333
$ws.Range("K74").Value = 1111.7931
$ws.Range("L74").Value = 2451.8333
$ws.Range("M74").Value = -237.7931000000001
$ws.Range("N74").Value = -4199.8333

$ws.Range("H77").Value = 1341.5143
$ws.Range("I77").Value = 1111.7931
$ws.Range("J77").Value = 2451.8333
$ws.Range("K77").Value = 5558.9655
$ws.Range("L77").Value = 12259.1665
$ws.Range("M77").Value = -1190.9655
$ws.Range("N77").Value = -20995.1665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1984.0625
$ws.Range("I20").Value = 1840
$ws.Range("J20").Value = 2416.25
$ws.Range("K20").Value = 1840
$ws.Range("L20").Value = 2416.25
$ws.Range("M20").Value = -1593
$ws.Range("N20").Value = -2910.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1816
$ws.Range("I31").Value = 999.3333
$ws.Range("J31").Value = 2428.5
$ws.Range("K31").Value = 999.3333
$ws.Range("L31").Value = 2428.5
$ws.Range("M31").Value = -704.3333
$ws.Range("N31").Value = -3018.5

$ws.Range("H33").Value = 1674.8334
$ws.Range("I33").Value = 1674.8334
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1674.8334
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1295.8334
$ws.Range("N33").ClearContents()

$ws.Range("H34").Value = 1816
$ws.Range("I34").Value = 999.3333
$ws.Range("J34").Value = 2428.5
$ws.Range("K34").Value = 999.3333
$ws.Range("L34").Value = 2428.5
$ws.Range("M34").Value = -797.3333
$ws.Range("N34").Value = -2832.5

$ws.Range("H107").Value = 384.94446
$ws.Range("I107").Value = 311.26666
$ws.Range("J107").Value = 753.3333
$ws.Range("K107").Value = 311.26666
$ws.Range("L107").Value = 753.3333
$ws.Range("M107").Value = 1608.73334
$ws.Range("N107").Value = -4593.3333

$ws.Range("H132").Value = 1677.4656
$ws.Range("I132").Value = 1147.7435
$ws.Range("J132").Value = 2764.7896
$ws.Range("K132").Value = 3443.2305
$ws.Range("L132").Value = 8294.3688
$ws.Range("M132").Value = -913.2305000000001
$ws.Range("N132").Value = -13354.3688

$ws.Range("H134").Value = 1785.3
$ws.Range("I134").Value = 1555.4073
$ws.Range("J134").Value = 3854.3333
$ws.Range("K134").Value = 4666.2219
$ws.Range("L134").Value = 11562.9999
$ws.Range("M134").Value = -2131.2219
$ws.Range("N134").Value = -16632.9999

$ws.Range("H135").Value = 45635.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 45635.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 45635.5
$ws.Range("N135").Value = -55775.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 863.625
$ws.Range("I32").Value = 798.6667
$ws.Range("J32").Value = 902.6
$ws.Range("K32").Value = 2396.0001
$ws.Range("L32").Value = 2707.8
$ws.Range("M32").Value = -2113.0001
$ws.Range("N32").Value = -3273.8

$ws.Range("H68").Value = 704.2857
$ws.Range("I68").Value = 646.6667
$ws.Range("J68").Value = 747.5
$ws.Range("K68").Value = 1940.0001
$ws.Range("L68").Value = 2242.5
$ws.Range("M68").Value = -1129.0001
$ws.Range("N68").Value = -3864.5

$ws.Range("H71").Value = 704.2857
$ws.Range("I71").Value = 646.6667
$ws.Range("J71").Value = 747.5
$ws.Range("K71").Value = 5820.0003
$ws.Range("L71").Value = 6727.5
$ws.Range("M71").Value = -1764.0003
$ws.Range("N71").Value = -14839.5

$ws.Range("H98").Value = 354.25
$ws.Range("I98").Value = 196.66667
$ws.Range("J98").Value = 448.8
$ws.Range("K98").Value = 590.00001
$ws.Range("L98").Value = 1346.4
$ws.Range("M98").Value = 907.99999
$ws.Range("N98").Value = -4342.4

$ws.Range("H107").Value = 726.95654
$ws.Range("I107").Value = 297.6
$ws.Range("J107").Value = 846.2222
$ws.Range("K107").Value = 892.8000000000001
$ws.Range("L107").Value = 2538.6666
$ws.Range("M107").Value = 1027.2
$ws.Range("N107").Value = -6378.6666

$ws.Range("H122").Value = 1009.1818
$ws.Range("I122").Value = 664
$ws.Range("J122").Value = 1296.8334
$ws.Range("K122").Value = 5976
$ws.Range("L122").Value = 11671.5006
$ws.Range("M122").Value = -3526
$ws.Range("N122").Value = -16571.5006

$ws.Range("H131").Value = 7834.758
$ws.Range("I131").Value = 472.7
$ws.Range("J131").Value = 8661.955
$ws.Range("K131").Value = 1418.1
$ws.Range("L131").Value = 25985.865
$ws.Range("M131").Value = 3621.9
$ws.Range("N131").Value = -36065.865

$ws.Range("H132").Value = 2416.3333
$ws.Range("I132").Value = 1749.5
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 15745.5
$ws.Range("L132").Value = 33750
$ws.Range("M132").Value = -13215.5
$ws.Range("N132").Value = -38810

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 9999
$ws.Range("I3").Value = 9999
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9999
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -9883
$ws.Range("N3").ClearContents()

$ws.Range("H7").Value = 5266706
$ws.Range("I7").Value = 5375000
$ws.Range("J7").Value = 5006800
$ws.Range("K7").Value = 5375000
$ws.Range("L7").Value = 5006800
$ws.Range("M7").Value = -5374888
$ws.Range("N7").Value = -5007024

$ws.Range("H8").Value = 5266706
$ws.Range("I8").Value = 5375000
$ws.Range("J8").Value = 5006800
$ws.Range("K8").Value = 5375000
$ws.Range("L8").Value = 5006800
$ws.Range("M8").Value = -5374861
$ws.Range("N8").Value = -5007078

$ws.Range("H70").Value = 24625
$ws.Range("I70").Value = 31333.334
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 31333.334
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -31063.334
$ws.Range("N70").Value = -5040

$ws.Range("H73").Value = 24625
$ws.Range("I73").Value = 31333.334
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 31333.334
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -30397.334
$ws.Range("N73").Value = -6372

$ws.Range("H102").Value = 3475.3
$ws.Range("I102").Value = 3656.625
$ws.Range("J102").Value = 2750
$ws.Range("K102").Value = 3656.625
$ws.Range("L102").Value = 2750
$ws.Range("M102").Value = -2034.625
$ws.Range("N102").Value = -5994

$ws.Range("H122").Value = 1387.4062
$ws.Range("I122").Value = 1125
$ws.Range("J122").Value = 1684.8
$ws.Range("K122").Value = 3375
$ws.Range("L122").Value = 5054.4
$ws.Range("M122").Value = -925
$ws.Range("N122").Value = -9954.4

$ws.Range("H133").Value = 48000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 48000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 48000
$ws.Range("N133").Value = -58120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4110.1816
$ws.Range("I40").Value = 1676.5
$ws.Range("J40").Value = 10600
$ws.Range("K40").Value = 1676.5
$ws.Range("L40").Value = 10600
$ws.Range("M40").Value = -1540.5
$ws.Range("N40").Value = -10872

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H132").Value = 1687.0227
$ws.Range("I132").Value = 1473.0952
$ws.Range("J132").Value = 1882.3478
$ws.Range("K132").Value = 4419.2856
$ws.Range("L132").Value = 5647.0434
$ws.Range("M132").Value = -1889.2856
$ws.Range("N132").Value = -10707.0434

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H138").Value = 88888
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 88888
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 88888
$ws.Range("N138").Value = -99168

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 18884
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 18884
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 18884
$ws.Range("N112").Value = -21838

$ws.Range("H132").Value = 1142.0488
$ws.Range("I132").Value = 832.9677
$ws.Range("J132").Value = 2100.2
$ws.Range("K132").Value = 2498.9031
$ws.Range("L132").Value = 6300.599999999999
$ws.Range("M132").Value = 31.09690000000001
$ws.Range("N132").Value = -11360.6
